$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell H74: 356 -> 357
$ws.Range("H74").Value = 357

# Add new row 75 for period 01-04-2021
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = -375
$ws.Range("C75").Value = -417
$ws.Range("D75").Value = -2
$ws.Range("E75").Value = -415
$ws.Range("F75").Value = 43
$ws.Range("G75").Value = 178
$ws.Range("H75").Value = -135
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = -2
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = -375
$ws.Range("O75").Value = -372
$ws.Range("P75").Value = -372
$ws.Range("Q75").Value = -3
